$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 30, shifting existing rows 30-51 down to 31-52
$ws.Rows.Item(30).Insert()

# Populate the newly inserted row 30 with the new weekly record
$ws.Cells.Item(30, 1).Value = 3
$ws.Cells.Item(30, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(30, 3).Value = "Coquimbo"
$ws.Cells.Item(30, 4).Value = 44893
$ws.Cells.Item(30, 5).Value = 5
$ws.Cells.Item(30, 6).Value = 300000000
$ws.Cells.Item(30, 7).Value = "Espárragos"
$ws.Cells.Item(30, 8).Value = "Verde"
$ws.Cells.Item(30, 9).Value = "Primera"
$ws.Cells.Item(30, 10).Value = 1230
$ws.Cells.Item(30, 11).Value = 1500
$ws.Cells.Item(30, 12).Value = 1600
$ws.Cells.Item(30, 13).Value = 1553
$ws.Cells.Item(30, 14).Value = "`$/kilo"
$ws.Cells.Item(30, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(30, 16).Value = 1553
$ws.Cells.Item(30, 17).Value = 1
$ws.Cells.Item(30, 18).Value = "Hortaliza"
